$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy cell formatting (styles) from column N into new columns O, P, Q ---
# Header row (row 1) uses the same style as N1
$ws.Range("N1").Copy()
$ws.Range("O1:Q1").PasteSpecial(-4122)

# Data rows (2-21): row 2 and rows 3-21 each carry their own style,
# which PasteSpecial(xlPasteFormats) preserves per source row.
$ws.Range("N2:N21").Copy()
$ws.Range("O2:Q21").PasteSpecial(-4122)

# --- New header labels ---
$ws.Range("O1").Value = "solarDeployed"
$ws.Range("P1").Value = "superchargerStationsAccumulated"
$ws.Range("Q1").Value = "superchargerConnectorsAccumulated"

# --- New data values ---
$ws.Range("O2").Value = 47
$ws.Range("P2").Value = 1490
$ws.Range("Q2").Value = 12767
$ws.Range("O3").Value = 29
$ws.Range("P3").Value = 1587
$ws.Range("Q3").Value = 13881
$ws.Range("O4").Value = 43
$ws.Range("P4").Value = 1653
$ws.Range("Q4").Value = 14658
$ws.Range("O5").Value = 54
$ws.Range("P5").Value = 1821
$ws.Range("Q5").Value = 16104
$ws.Range("O6").Value = 35
$ws.Range("P6").Value = 1917
$ws.Range("Q6").Value = 17007
$ws.Range("O7").Value = 27
$ws.Range("P7").Value = 2035
$ws.Range("Q7").Value = 18100
$ws.Range("O8").Value = 57
$ws.Range("P8").Value = 2181
$ws.Range("Q8").Value = 19437
$ws.Range("O9").Value = 86
$ws.Range("P9").Value = 2564
$ws.Range("Q9").Value = 23277
$ws.Range("O10").Value = 92
$ws.Range("P10").Value = 2699
$ws.Range("Q10").Value = 24515
$ws.Range("O11").Value = 85
$ws.Range("P11").Value = 2966
$ws.Range("Q11").Value = 26900
$ws.Range("O12").Value = 83
$ws.Range("P12").Value = 3254
$ws.Range("Q12").Value = 29281
$ws.Range("O13").Value = 85
$ws.Range("P13").Value = 3476
$ws.Range("Q13").Value = 31498
$ws.Range("O14").Value = 48
$ws.Range("P14").Value = 3724
$ws.Range("Q14").Value = 33657
$ws.Range("O15").Value = 106
$ws.Range("P15").Value = 3971
$ws.Range("Q15").Value = 36165
$ws.Range("O16").Value = 94
$ws.Range("P16").Value = 4283
$ws.Range("Q16").Value = 38883
$ws.Range("O17").Value = 100
$ws.Range("P17").Value = 4678
$ws.Range("Q17").Value = 42419
$ws.Range("O18").Value = 67
$ws.Range("P18").Value = 4947
$ws.Range("Q18").Value = 45169
$ws.Range("O19").Value = 66
$ws.Range("P19").Value = 5265
$ws.Range("Q19").Value = 48082
$ws.Range("O20").Value = 49
$ws.Range("P20").Value = 5595
$ws.Range("Q20").Value = 51105
$ws.Range("O21").Value = 41
$ws.Range("P21").Value = 5952
$ws.Range("Q21").Value = 54892
